$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns in this sheet are plain text
# (e.g. "557.34", "  +3.02%  "), not real numbers/percentages. Assigning a
# numeric-looking string straight to .Value lets Excel auto-coerce it into a
# Number, which would lose the original text formatting (e.g. "1.00" -> 1).
# Set-CellText avoids that by temporarily forcing the cell to Text format
# ("@") before writing the value, then restores whatever format the cell had.
function Set-CellText($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    $prevFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $prevFormat
}

Set-CellText $ws 'D2' '63.988.00'
Set-CellText $ws 'E2' '  +3.49%  '
Set-CellText $ws 'D3' '3.054.08'
Set-CellText $ws 'E3' '  +2.27%  '
Set-CellText $ws 'E4' '  -0.07%  '
Set-CellText $ws 'D5' '557.34'
Set-CellText $ws 'E5' '  +3.02%  '
Set-CellText $ws 'D6' '142.61'
Set-CellText $ws 'E6' '  +5.61%  '
Set-CellText $ws 'E7' '  +0.05%  '
Set-CellText $ws 'D8' '3.048.68'
Set-CellText $ws 'E8' '  +2.14%  '
Set-CellText $ws 'D9' '0.512'
Set-CellText $ws 'E9' '  +5.57%  '
Set-CellText $ws 'D10' '0.155'
Set-CellText $ws 'E10' '  +6.57%  '
Set-CellText $ws 'D11' '6.05'
Set-CellText $ws 'E11' '  -9.66%  '
Set-CellText $ws 'E12' '  +7.68%  '
Set-CellText $ws 'E13' '  +6.39%  '
Set-CellText $ws 'D14' '35.04'
Set-CellText $ws 'E14' '  +4.44%  '
Set-CellText $ws 'D15' '3.540.16'
Set-CellText $ws 'E15' '  +2.97%  '
Set-CellText $ws 'D16' '64.020.56'
Set-CellText $ws 'E16' '  +3.49%  '
Set-CellText $ws 'D17' '3.055.93'
Set-CellText $ws 'E17' '  +2.55%  '
Set-CellText $ws 'E18' '  +1.81%  '
Set-CellText $ws 'E19' '  +2.89%  '
Set-CellText $ws 'D20' '475.89'
Set-CellText $ws 'E20' '  +2.60%  '
Set-CellText $ws 'E21' '  +4.84%  '
Set-CellText $ws 'E22' '  +4.11%  '
Set-CellText $ws 'D23' '7.56'
Set-CellText $ws 'E23' '  +6.06%  '
Set-CellText $ws 'D24' '14.19'
Set-CellText $ws 'E24' '  +13.81%  '
Set-CellText $ws 'D25' '81.70'
Set-CellText $ws 'E25' '  +3.25%  '
Set-CellText $ws 'E26' '  -0.12%  '
Set-CellText $ws 'E27' '  +3.14%  '
Set-CellText $ws 'D28' '7.90'
Set-CellText $ws 'E28' '  +5.20%  '
Set-CellText $ws 'E29' '  +2.39%  '
Set-CellText $ws 'D30' '1.00'
Set-CellText $ws 'E30' '  +0.07%  '
Set-CellText $ws 'E31' '  +4.28%  '
Set-CellText $ws 'E32' '  +1.80%  '
Set-CellText $ws 'E33' '  +5.98%  '
Set-CellText $ws 'E34' '  +2.61%  '
Set-CellText $ws 'D35' '6.18'
Set-CellText $ws 'D36' '54.64'
Set-CellText $ws 'E36' '  +1.31%  '
Set-CellText $ws 'D37' '0.0404'
Set-CellText $ws 'E37' '  +5.06%  '
Set-CellText $ws 'D38' '441.36'
Set-CellText $ws 'E38' '  -1.31%  '
Set-CellText $ws 'D39' '0.0804'
Set-CellText $ws 'E39' '  +0.50%  '
Set-CellText $ws 'D40' '2.82'
Set-CellText $ws 'E40' '  +16.03%  '
Set-CellText $ws 'D41' '2.964.00'
Set-CellText $ws 'E41' '  +1.15%  '
Set-CellText $ws 'E42' '  +3.16%  '
Set-CellText $ws 'E43' '  -0.94%  '
Set-CellText $ws 'D44' '27.56'
Set-CellText $ws 'E44' '  +4.67%  '
Set-CellText $ws 'E45' '  +6.19%  '
Set-CellText $ws 'E46' '  +8.89%  '
Set-CellText $ws 'E48' '  +4.86%  '
Set-CellText $ws 'B49' 'PEPE'
Set-CellText $ws 'C49' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText $ws 'D49' '0.0₃0513'
Set-CellText $ws 'E49' '  +5.83%  '
Set-CellText $ws 'B50' 'Monero'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText $ws 'D50' '117.12'
Set-CellText $ws 'E50' '  +2.23%  '
Set-CellText $ws 'D51' '2.06'
Set-CellText $ws 'E51' '  +3.94%  '
